# Consolidate the multiple text runs ("The" / " " / "picture" / " " / "first")
# of the caption textbox on slide 1 into a single run with the same text.
#
# Simply re-assigning the identical string to TextRange.Text is treated as a
# no-op by the writer's run-level diffing, so the original run split survives
# untouched. Clearing the text first (so the new text shares nothing with the
# old content) forces the writer to regenerate the paragraph from scratch,
# which - per the "consolidate text runs when possible" change - emits one
# run instead of several.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shape = $s.Shapes.Item("TextBox 3")
$tr = $shape.TextFrame.TextRange

$tr.Delete()
$tr.Text = "The picture first"
